$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 12:05"

# Update country-stat cells per the refreshed COVID data snapshot
$ws.Range("B22").Value = 58517
$ws.Range("C22").Value = 136
$ws.Range("D22").Value = 15919
$ws.Range("E22").Value = 33112
$ws.Range("G22").Value = 19
$ws.Range("H22").Value = 9486
$ws.Range("B36").Value = 26940
$ws.Range("C36").Value = 467
$ws.Range("D36").Value = 7637
$ws.Range("E36").Value = 17662
$ws.Range("G36").Value = 28
$ws.Range("H36").Value = 1641
$ws.Range("B41").Value = 19398
$ws.Range("C41").Value = 141
$ws.Range("D41").Value = 13426
$ws.Range("E41").Value = 4702
$ws.Range("B42").Value = 18638
$ws.Range("C42").Value = 552
$ws.Range("D42").Value = 3979
$ws.Range("E42").Value = 13699
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 960
$ws.Range("B64").Value = 7857
$ws.Range("C64").Value = 38
$ws.Range("D64").Value = 6404
$ws.Range("E64").Value = 1338
$ws.Range("B65").Value = 7819
$ws.Range("C65").Value = 12
$ws.Range("D65").Value = 5754
$ws.Range("E65").Value = 1860
$ws.Range("B67").Value = 6885
$ws.Range("C67").Value = 26
$ws.Range("E67").Value = 1065
$ws.Range("E106").Value = 6
$ws.Range("G106").Value = 1
$ws.Range("H106").Value = 109
$ws.Range("B112").Value = 1143
$ws.Range("C112").Value = 6
$ws.Range("D112").Value = 877
$ws.Range("E112").Value = 233
$ws.Range("B113").Value = 1088
$ws.Range("C113").Value = 3
$ws.Range("E113").Value = 47
$ws.Range("B114").Value = 1084
$ws.Range("C114").Value = 7
$ws.Range("D114").Value = 964
$ws.Range("E114").Value = 72
$ws.Range("A142").Value = "Uganda"
$ws.Range("B142").Value = 457
$ws.Range("C142").Value = 40
$ws.Range("D142").Value = 72
$ws.Range("E142").Value = 385
$ws.Range("H142").Value = 0
$ws.Range("A143").Value = "Estado de Palestina"
$ws.Range("B143").Value = 449
$ws.Range("D143").Value = 372
$ws.Range("E143").Value = 74
$ws.Range("H143").Value = 3
$ws.Range("A144").Value = "Taiwan"
$ws.Range("B144").Value = 443
$ws.Range("C144").Value = 1
$ws.Range("D144").Value = 427
$ws.Range("E144").Value = 9
$ws.Range("H144").Value = 7
$ws.Range("A145").Value = "Togo"
$ws.Range("B145").Value = 442
$ws.Range("D145").Value = 211
$ws.Range("E145").Value = 218
$ws.Range("H145").Value = 13
$ws.Range("A146").Value = "Cabo Verde"
$ws.Range("B146").Value = 435
$ws.Range("D146").Value = 193
$ws.Range("E146").Value = 238
$ws.Range("H146").Value = 4
$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("D202").Value = 18
$ws.Range("H202").Value = 0
$ws.Range("B215").Value = 7
$ws.Range("D215").Value = 7
